{"js": "// Replace the multiplication-problem text in each table cell with the new\n// problem/answer pairs, per the commit's regenerated worksheet values.\n// Each old value is unique in the document, so a targeted search + replace\n// on the matched range keeps the existing run/paragraph formatting intact.\nconst replacements = [\n  [\"90\u00d743=3870\", \"70\u00d723=1610\"],\n  [\"93\u00d794=8742\", \"89\u00d755=4895\"],\n  [\"71\u00d745=3195\", \"37\u00d799=3663\"],\n  [\"43\u00d789=3827\", \"73\u00d721=1533\"],\n  [\"52\u00d713=676\", \"56\u00d733=1848\"],\n  [\"16\u00d751=816\", \"21\u00d758=1218\"],\n  [\"43\u00d717=731\", \"20\u00d789=1780\"],\n  [\"53\u00d730=1590\", \"83\u00d721=1743\"],\n  [\"28\u00d752=1456\", \"61\u00d737=2257\"],\n  [\"24\u00d753=1272\", \"28\u00d751=1428\"],\n  [\"13\u00d737=481\", \"82\u00d736=2952\"],\n  [\"95\u00d782=7790\", \"31\u00d780=2480\"],\n  [\"21\u00d798=2058\", \"39\u00d785=3315\"],\n  [\"75\u00d793=6975\", \"74\u00d790=6660\"],\n  [\"60\u00d712=720\", \"68\u00d784=5712\"],\n  [\"20\u00d762=1240\", \"82\u00d778=6396\"],\n  [\"99\u00d784=8316\", \"89\u00d734=3026\"],\n  [\"18\u00d718=324\", \"75\u00d781=6075\"],\n  [\"74\u00d791=6734\", \"20\u00d738=760\"],\n  [\"17\u00d776=1292\", \"64\u00d743=2752\"],\n  [\"62\u00d758=3596\", \"64\u00d718=1152\"],\n  [\"55\u00d724=1320\", \"12\u00d733=396\"],\n  [\"11\u00d743=473\", \"84\u00d769=5796\"],\n  [\"47\u00d759=2773\", \"12\u00d770=840\"],\n  [\"77\u00d713=1001\", \"86\u00d795=8170\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication-problem text in each table cell with the new\n# problem/answer pairs, per the commit's regenerated worksheet values.\n# Each old value is unique in the document, so Find/Replace across the\n# whole document body safely retargets only the intended cell each time,\n# leaving paragraph/run formatting untouched.\n$pairs = @(\n    @{old=\"90\u00d743=3870\"; new=\"70\u00d723=1610\"},\n    @{old=\"93\u00d794=8742\"; new=\"89\u00d755=4895\"},\n    @{old=\"71\u00d745=3195\"; new=\"37\u00d799=3663\"},\n    @{old=\"43\u00d789=3827\"; new=\"73\u00d721=1533\"},\n    @{old=\"52\u00d713=676\";  new=\"56\u00d733=1848\"},\n    @{old=\"16\u00d751=816\";  new=\"21\u00d758=1218\"},\n    @{old=\"43\u00d717=731\";  new=\"20\u00d789=1780\"},\n    @{old=\"53\u00d730=1590\"; new=\"83\u00d721=1743\"},\n    @{old=\"28\u00d752=1456\"; new=\"61\u00d737=2257\"},\n    @{old=\"24\u00d753=1272\"; new=\"28\u00d751=1428\"},\n    @{old=\"13\u00d737=481\";  new=\"82\u00d736=2952\"},\n    @{old=\"95\u00d782=7790\"; new=\"31\u00d780=2480\"},\n    @{old=\"21\u00d798=2058\"; new=\"39\u00d785=3315\"},\n    @{old=\"75\u00d793=6975\"; new=\"74\u00d790=6660\"},\n    @{old=\"60\u00d712=720\";  new=\"68\u00d784=5712\"},\n    @{old=\"20\u00d762=1240\"; new=\"82\u00d778=6396\"},\n    @{old=\"99\u00d784=8316\"; new=\"89\u00d734=3026\"},\n    @{old=\"18\u00d718=324\";  new=\"75\u00d781=6075\"},\n    @{old=\"74\u00d791=6734\"; new=\"20\u00d738=760\"},\n    @{old=\"17\u00d776=1292\"; new=\"64\u00d743=2752\"},\n    @{old=\"62\u00d758=3596\"; new=\"64\u00d718=1152\"},\n    @{old=\"55\u00d724=1320\"; new=\"12\u00d733=396\"},\n    @{old=\"11\u00d743=473\";  new=\"84\u00d769=5796\"},\n    @{old=\"47\u00d759=2773\"; new=\"12\u00d770=840\"},\n    @{old=\"77\u00d713=1001\"; new=\"86\u00d795=8170\"}\n)\n\n$d = $word.ActiveDocument\n\nforeach ($p in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $p.old\n    $find.Replacement.Text = $p.new\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
